# The capital call "Call Date" (D2) and "Due Date" (E2) values are updated:
#  - D2 becomes an explicit date value formatted as a short date (instead of
#    being unformatted with the sheet's default numeric formatting).
#  - E2 stops being computed via the formula "=D2+7" and instead holds a
#    plain, hard-coded date value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCall")

# Give D2 a short-date number format while keeping the normal body font
# (copy the formatting from a neighboring default-styled cell so the same
# font is reused rather than a new one being created), then set its value.
$ws.Range("I2").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Value = 44856

# E2 is no longer a formula; it now just holds a fixed date value.
$ws.Range("E2").Value = 44870
